$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. The old last paragraph ("Mostrado o conceito de f" + bookmark + "actory...")
#    already reads correctly as one sentence once the two runs are
#    concatenated. Force Word to re-flow / merge those two runs into a
#    single run by doing a no-op Find&Replace over text that spans the
#    run boundary. As a side effect this also removes the _GoBack bookmark
#    that used to sit at the old edit point, between the two runs - it will
#    be re-added at the very end of the document, after the newly appended
#    text, which is where Word leaves it after the newest edit.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("de factory", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "de factory", 2) | Out-Null

# ---------------------------------------------------------------------------
# Helper text to append, modelled after the existing "Nª aula:" sections
# already present in the document.
# ---------------------------------------------------------------------------

# Blank separator paragraph, then "10ª aula:" heading, then body text.
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()

$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$heading = $d.Paragraphs.Last
$heading.Range.Font.Bold = $true
$heading.Range.Font.BoldBi = $true
$heading.Range.Font.Size = 14
$heading.Range.Font.SizeBi = 14
$heading.Range.InsertAfter("10ª aula:")

$heading.Range.InsertParagraphAfter()
$body = $d.Paragraphs.Last
$body.Range.InsertAfter("Inicio do projeto de criador de formulários. Criada classe do input e definido atributos type e required e o método construtor.")

# Blank separator paragraph, then "11ª aula:" heading, then body text.
$body.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$heading = $d.Paragraphs.Last
$heading.Range.Font.Bold = $true
$heading.Range.Font.BoldBi = $true
$heading.Range.Font.Size = 14
$heading.Range.Font.SizeBi = 14
$heading.Range.InsertAfter("11ª aula:")

$heading.Range.InsertParagraphAfter()
$body = $d.Paragraphs.Last
$body.Range.InsertAfter("Criação das classes de botão e formulário. Botão é uma classe filha de input, herdando seus valores e a classe formulário recebe os objetos de input e botão.  ")

# Blank separator paragraph, then "12ª aula:" heading, then body text.
$body.Range.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
$heading = $d.Paragraphs.Last
$heading.Range.Font.Bold = $true
$heading.Range.Font.BoldBi = $true
$heading.Range.Font.Size = 14
$heading.Range.Font.SizeBi = 14
$heading.Range.InsertAfter("12ª aula:")

$heading.Range.InsertParagraphAfter()
$body = $d.Paragraphs.Last
$body.Range.InsertAfter("Finalização do projeto. Criado os métodos de render, exibindo os elementos na tela.")

# ---------------------------------------------------------------------------
# 3. Re-add the _GoBack bookmark at the very end of the document, marking
#    the location of this latest edit - same convention Word uses.
# ---------------------------------------------------------------------------
$endRange = $d.Paragraphs.Last.Range
$endOfDoc = $d.Range($endRange.End - 1, $endRange.End - 1)
$d.Bookmarks.Add("_GoBack", $endOfDoc)
